$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking)
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total)
$ws.Range("B12").Value = 95
$ws.Range("C12").Value = -10.8
$ws.Range("E12").Value = "84.2/140"
